# Trade #40 closed at 2026-02-16 22:56:06 - base_strategy DOWN +0.000%
#
# Appends a new trade row (row 41) to both the "All Trades" and
# "base_strategy" worksheets, mirroring the existing row 40 layout.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 41

    $ws.Range("A$row").Value = 40

    # Date / Time columns hold plain text that looks like a date/time -
    # force Text format first so Excel doesn't auto-convert them to
    # date/time serial numbers.
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = "2026-02-16"

    $ws.Range("C$row").NumberFormat = "@"
    $ws.Range("C$row").Value = "22:56:06"

    $ws.Range("D$row").Value = "base_strategy"
    $ws.Range("E$row").Value = "DOWN"
    $ws.Range("F$row").Value = 49.999998

    # G (Exit Price) stays blank - trade is still OPEN.
    $ws.Range("H$row").Value = "OPEN"
    $ws.Range("I$row").Value = 0
    $ws.Range("J$row").Value = 0
    $ws.Range("K$row").Value = 100
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0.6
    $ws.Range("O$row").Value = "Normal spread capture: 19600 bps"

    # P (Exit Reason) stays blank - trade is still OPEN.
    $ws.Range("Q$row").Value = 0
}
